# issue #5: stock data from json to db
#
# The "股票" (stock) worksheet (sheet4) gains three new columns:
#   - "category"    inserted right after "property_category" (pushes
#                    date / legislator_name / legislator_id one column right)
#   - "source_file"  appended after "legislator_id"
#   - "index"        appended after "source_file" (mirrors the row's id
#                    that already lives in column A)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

$firstDataRow = 2
$lastDataRow = 13

# Column I currently holds "date" (property_category=H, date=I, legislator_name=J,
# legislator_id=K). Insert a new blank column at I so everything from the old
# I onward shifts right by one; the Insert also carries over each row's
# existing cell style (header style for row 1, data style for the rest).
$ws.Columns.Item(9).Insert() | Out-Null

# New column I = "category"
$ws.Range("I1").Value = "category"
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $ws.Range("I$r").Value = "normal"
}

# Append two more columns: M = "source_file", N = "index".
# Copy formatting from column L (legislator_id) so the new cells pick up the
# same header/data styling already used across the sheet.
$ws.Range("L1:L$lastDataRow").Copy() | Out-Null
$ws.Range("M1:M$lastDataRow").PasteSpecial(-4122) | Out-Null
$ws.Range("L1:L$lastDataRow").Copy() | Out-Null
$ws.Range("N1:N$lastDataRow").PasteSpecial(-4122) | Out-Null

$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $ws.Range("M$r").Value = "tmp25ce1"
    $rowIndex = $ws.Range("A$r").Value()
    $ws.Range("N$r").Value = $rowIndex
}
